$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new column at N ("下学期艺考赋分"), pushing the existing
# N ("常规教学薄成绩") and O ("备注") columns one slot to the right.
$ws.Columns("N").EntireColumn.Insert()

# Header for the new column.
$ws.Range("N1").Value = "下学期艺考赋分"

# New column keeps the same centered look as the rest of the header/data
# cells (style index 1 in the original workbook).
$ws.Range("N1:N9").HorizontalAlignment = -4108

# Fill the new column's data rows with the value 9.
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 14).Value = 9
}

# Restore the explicit column width for the new N column (14 characters);
# the shifted-right former N column (now O) keeps its original 17.5 width
# automatically.
$ws.Columns("N").ColumnWidth = 13.29

# Matches the recorded selection after the edit.
$ws.Range("N12").Select() | Out-Null
